# The post "「かきとりのちから」" (row 674) was removed from the sheet.
# Deleting the entire row shifts every row below it up by one, which
# matches the target state (dimension shrinks from A1:C744 to A1:C743
# and all subsequent rows are renumbered accordingly, with no other
# content changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("674").Delete()
